# Apply the two-digit-mul.docx content updates:
#  - update the date heading
#  - replace each multiplication prompt in the table with its new value

$d = $word.ActiveDocument

$replacements = @(
    @{ Old = "2023-09-02 Saturday"; New = "2023-09-03 Sunday" },
    @{ Old = "77×17=";  New = "76×54=" },
    @{ Old = "92×84=";  New = "13×80=" },
    @{ Old = "50×49=";  New = "46×15=" },
    @{ Old = "36×21=";  New = "61×94=" },
    @{ Old = "47×54=";  New = "93×19=" },
    @{ Old = "39×26=";  New = "52×11=" },
    @{ Old = "48×78=";  New = "80×36=" },
    @{ Old = "68×75=";  New = "88×73=" },
    @{ Old = "47×70=";  New = "14×89=" },
    @{ Old = "92×64=";  New = "44×41=" },
    @{ Old = "58×68=";  New = "88×73=" },
    @{ Old = "31×76=";  New = "87×75=" },
    @{ Old = "51×37=";  New = "82×33=" },
    @{ Old = "82×85=";  New = "36×33=" },
    @{ Old = "22×26=";  New = "35×21=" },
    @{ Old = "94×42=";  New = "55×41=" },
    @{ Old = "54×95=";  New = "61×30=" },
    @{ Old = "65×31=";  New = "88×90=" },
    @{ Old = "60×81=";  New = "32×92=" },
    @{ Old = "28×54=";  New = "27×17=" },
    @{ Old = "93×34=";  New = "30×97=" },
    @{ Old = "88×32=";  New = "28×43=" },
    @{ Old = "78×69=";  New = "89×92=" },
    @{ Old = "35×54=";  New = "13×20=" },
    @{ Old = "36×56=";  New = "54×19=" }
)

foreach ($r in $replacements) {
    $find = $d.Content.Find
    $find.ClearFormatting()
    $find.Execute($r.Old, $true, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)
}
